$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.960.99"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -1.28%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'3.166.48"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -4.57%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'591.28"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -2.14%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'134.62"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -5.56%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  +0.12%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'3.162.53"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -4.66%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  -1.23%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  -5.49%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  -5.15%  "
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = "'  -3.16%  "
$ws.Range("E12").ClearFormats()
$ws.Range("E13").Value = "'  -4.28%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'34.63"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -0.69%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'3.688.53"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -4.61%  "
$ws.Range("E15").ClearFormats()
$ws.Range("E16").Value = "'  -1.79%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'3.170.98"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -4.55%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'62.908.92"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -1.48%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  -4.54%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'460.24"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -4.23%  "
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'  -1.80%  "
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "'  -5.18%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'7.60"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -6.92%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'13.37"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").Value = "'83.52"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -1.48%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'  -0.05%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'  +0.04%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'2.67"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -3.72%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'7.68"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -6.41%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  -6.82%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  -6.27%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'27.13"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -6.16%  "
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = "'  -2.98%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'2.37"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -6.60%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'1.03"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -6.46%  "
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = "'  -4.10%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'51.18"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -2.26%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.0₃0702"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -5.18%  "
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = "'  -3.02%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'401.13"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -7.42%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E42").Value = "'  -4.95%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'  -8.10%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'2.791.70"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -10.10%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.250"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -5.41%  "
$ws.Range("E45").ClearFormats()
$ws.Range("E47").Value = "'  -5.28%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'123.69"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -0.13%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'25.30"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -4.02%  "
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'  -2.34%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'34.28"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -8.22%  "
$ws.Range("E51").ClearFormats()
